$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D..AQ for data rows 2 and 3 (headers are in row 1).
# Column letter => new value
$values = [ordered]@{
    "D"  = 0.033
    "E"  = 0.036
    "G"  = 0.2302457466918715
    "H"  = 0.2302457466918715
    "I"  = 0.2022684310018903
    "J"  = 0.1348456206679269
    "K"  = 32.1
    "L"  = 0.1213610586011342
    "M"  = 39.7
    "N"  = 0.05209290119406902
    "O"  = 1.236760124610592
    "P"  = 39.7
    "Q"  = 0.05209290119406902
    "R"  = 1.236760124610592
    "S"  = 0
    "T"  = 0
    "U"  = 4.18
    "V"  = 0.005484844508594672
    "W"  = 0.2316017316017316
    "X"  = 0.06193192762843093
    "Y"  = 0.1696698039733007
    "Z"  = 2.10287804102401
    "AA" = 0.2835638946308369
    "AB" = 0.0617373318700806
    "AC" = 0.2218265627607563
    "AD" = 4.46
    "AE" = 0
    "AF" = 4.46
    "AG" = 0.2800000000000002
    "AH" = 0.00581820079315383
    "AI" = 0.03141730064806988
    "AJ" = 0.0003672709147669145
    "AK" = 0.002032225286688926
    "AL" = 0
    "AM" = -1.89
    "AN" = 0.08050541516245488
    "AP" = 0.005054151624548741
    "AQ" = -28.30687830687831
}

foreach ($row in @(2, 3)) {
    foreach ($col in $values.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $values[$col]
    }
}
